$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "日期：2018.10.10 第六周周一" block (rows 31-39), mirroring the
# existing "日期：2018.10.10 第六周周一" block at rows 20-28 (same header date,
# same team roster / task text), but with column C turned into a 100% done
# indicator for the first three members.
# ---------------------------------------------------------------------------

# 1) Merge the two banner rows FIRST (A31:D31 / A38:D39), matching how the
#    source block's equivalent rows are merged - doing this before copying
#    formats keeps a single uniform style across the merged row instead of
#    Excel splitting per-cell borders for a "new" merge.
$ws.Range("A31:D31").Merge()
$ws.Range("A38:D39").Merge()

# 2) Clone all formatting (borders/fonts/number formats/etc.) for the whole
#    A31:D39 block from the analogous A20:D28 block in one shot.
$ws.Range("A20:D28").Copy()
$ws.Range("A31").PasteSpecial(-4122)

# 3) Row 31 - the merged date banner (mirrors row 20).
$ws.Cells.Item(31, 1).Value2 = $ws.Cells.Item(20, 1).Value2

# 4) Row 32 - the column headers (mirrors row 21): 组员/计划内容/完成情况/备注.
$ws.Cells.Item(32, 1).Value2 = $ws.Cells.Item(21, 1).Value2
$ws.Cells.Item(32, 2).Value2 = $ws.Cells.Item(21, 2).Value2
$ws.Cells.Item(32, 3).Value2 = $ws.Cells.Item(21, 3).Value2
$ws.Cells.Item(32, 4).Value2 = $ws.Cells.Item(21, 4).Value2

# 5) Rows 33-37 - member + task columns mirror rows 22-26 exactly.
for ($i = 0; $i -le 4; $i++) {
    $srcRow = 22 + $i
    $dstRow = 33 + $i
    $ws.Cells.Item($dstRow, 1).Value2 = $ws.Cells.Item($srcRow, 1).Value2
    $ws.Cells.Item($dstRow, 2).Value2 = $ws.Cells.Item($srcRow, 2).Value2
}

# 6) Column C for rows 33-37 becomes a 100% "done" indicator instead of text.
#    Reuse the number-format/border combos that already exist on the sheet:
#    C22's style for rows 33-35, C3's style for rows 36-37 (same combos the
#    original workbook uses elsewhere for percentage cells).
$ws.Range("C22").Copy()
$ws.Range("C33:C35").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C36:C37").PasteSpecial(-4122)
$ws.Range("C33:C37").Value2 = 1

# 7) Row 38-39 - the merged "总结：" summary block (mirrors rows 27-28).
$ws.Cells.Item(38, 1).Value2 = $ws.Cells.Item(27, 1).Value2

# 8) Move the selection down to the newly added summary block and scroll the
#    viewport so the new rows are visible (mirrors the saved view state).
$ws.Range("A38:D39").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

Write-Output "done"
